$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: swap the given column values (by column index) between two rows.
function Swap-RowValues {
    param($row1, $row2, $cols)
    foreach ($col in $cols) {
        $v1 = $ws.Cells.Item($row1, $col).Value()
        $v2 = $ws.Cells.Item($row2, $col).Value()
        $ws.Cells.Item($row1, $col).Value = $v2
        $ws.Cells.Item($row2, $col).Value = $v1
    }
}

# Helper: rotate the given column values (by column index) across a list of rows,
# so that each row receives the values previously held by the row before it
# (cyclically) - i.e. rows[i] <- old rows[i-1], with wraparound.
function Rotate-RowValues {
    param($rows, $cols)

    $old = @{}
    foreach ($r in $rows) {
        foreach ($c in $cols) {
            $old["$r-$c"] = $ws.Cells.Item($r, $c).Value()
        }
    }

    for ($i = 0; $i -lt $rows.Length; $i++) {
        $target = $rows[$i]
        $srcIndex = ($i - 1 + $rows.Length) % $rows.Length
        $source = $rows[$srcIndex]
        foreach ($c in $cols) {
            $ws.Cells.Item($target, $c).Value = $old["$source-$c"]
        }
    }
}

# Columns B, E, F, G correspond to indices 2, 5, 6, 7.
$colsBEFG = @(2,5,6,7)
# Columns B, D, E, F, G correspond to indices 2, 4, 5, 6, 7.
$colsBDEFG = @(2,4,5,6,7)

# Row pairs whose B/E/F/G values were swapped (stock batch rows exchanged).
$pairs = @(
    @(136,137),
    @(163,164),
    @(233,234),
    @(246,247),
    @(292,293),
    @(299,300),
    @(311,312),
    @(465,466),
    @(467,468),
    @(469,470),
    @(472,473),
    @(479,480)
)

foreach ($pair in $pairs) {
    Swap-RowValues $pair[0] $pair[1] $colsBEFG
}

# Row pair whose B/D/E/F/G values were swapped.
Swap-RowValues 420 421 $colsBDEFG

# Rows 294/295/296 rotate values among themselves (294 gets 296's old data,
# 295 gets 294's old data, 296 gets 295's old data).
Rotate-RowValues @(294,295,296) $colsBEFG
